$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.177.14"
$ws.Range("E2").Value = "  +3.10%  "
$ws.Range("D3").Value = "1.579.67"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.38%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.64%  "
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("D13").Value = "1.578.50"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").Value = "29.213.51"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "0.0₃0689"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("E26").Value = "  +4.75%  "
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0469"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D33").Value = "1.422.81"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("E37").Value = "  +6.26%  "
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("E40").Value = "  +3.54%  "
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "52.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +23.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.790"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0471"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "1.718.28"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.841"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("E51").Value = "  -0.63%  "
